$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Fill in the new time-log entry on row 22 (previously a blank Table1 row):
# Date 2017-11-19, Start 14:15 (0.59375), End 16:00 (0.66666666666666663).
# The Work Time column (D) is a shared formula (=ABS(C-B)) that recalculates
# automatically, and the table's totals row (D36) picks up the new total.
$ws.Range("A22").Value = 43058
$ws.Range("B22").Value = 0.59375
$ws.Range("C22").Value = 0.66666666666666663

# Scroll the view up and move the active selection to E22
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("E22").Select()
